$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header swap: BP1/BQ1 (average_doctor <-> average_doctor_old)
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Updated statistics values (Harvard case classification)
# Row 4
$ws.Range("E4").Value = 0.427
$ws.Range("F4").Value = 0.07099999999999999
$ws.Range("G4").Value = 0.267
$ws.Range("N4").Value = 0.437
$ws.Range("O4").Value = 0.064
$ws.Range("P4").Value = 0.254
$ws.Range("Q4").Value = 0.024
$ws.Range("R4").Value = 0.017
$ws.Range("S4").Value = 0.13
$ws.Range("W4").Value = 0.292
$ws.Range("X4").Value = 0.109
$ws.Range("Y4").Value = 0.331
$ws.Range("AI4").Value = 0.292
$ws.Range("AJ4").Value = 0.089
$ws.Range("AK4").Value = 0.298
$ws.Range("AU4").Value = 0.194
$ws.Range("AV4").Value = 0.028
$ws.Range("AW4").Value = 0.167
$ws.Range("BA4").Value = 1.985
$ws.Range("BB4").Value = 0.159
$ws.Range("BC4").Value = 0.399
$ws.Range("BG4").Value = 0.723
$ws.Range("BH4").Value = 0.143
$ws.Range("BI4").Value = 0.378
$ws.Range("BM4").Value = 0.71
$ws.Range("BN4").Value = 0.08
$ws.Range("BO4").Value = 0.282
$ws.Range("BP4").Value = 0.662
$ws.Range("BQ4").Value = 0.704
# Row 5
$ws.Range("E5").Value = 0.543
$ws.Range("F5").Value = 0.08500000000000001
$ws.Range("G5").Value = 0.291
$ws.Range("N5").Value = 0.738
$ws.Range("O5").Value = 0.077
$ws.Range("P5").Value = 0.278
$ws.Range("Q5").Value = 0.016
$ws.Range("R5").Value = 0.007
$ws.Range("S5").Value = 0.08400000000000001
$ws.Range("W5").Value = 0.282
$ws.Range("X5").Value = 0.11
$ws.Range("Y5").Value = 0.331
$ws.Range("AI5").Value = 0.314
$ws.Range("AJ5").Value = 0.099
$ws.Range("AK5").Value = 0.315
$ws.Range("AU5").Value = 0.377
$ws.Range("AV5").Value = 0.092
$ws.Range("AW5").Value = 0.303
$ws.Range("BA5").Value = 1.334
$ws.Range("BB5").Value = 0.081
$ws.Range("BC5").Value = 0.285
$ws.Range("BG5").Value = 0.393
$ws.Range("BI5").Value = 0.227
$ws.Range("BM5").Value = 0.552
$ws.Range("BN5").Value = 0.065
$ws.Range("BO5").Value = 0.255
$ws.Range("BP5").Value = 0.445
$ws.Range("BQ5").Value = 0.455
# Row 6
$ws.Range("E6").Value = 0.478
$ws.Range("N6").Value = 0.549
$ws.Range("Q6").Value = 0.019
$ws.Range("W6").Value = 0.287
$ws.Range("AI6").Value = 0.303
$ws.Range("AU6").Value = 0.256
$ws.Range("BA6").Value = 1.586
$ws.Range("BG6").Value = 0.509
$ws.Range("BM6").Value = 0.621
$ws.Range("BP6").Value = 0.529
$ws.Range("BQ6").Value = 0.549
# Row 7
$ws.Range("E7").Value = 0.515
$ws.Range("N7").Value = 0.649
$ws.Range("Q7").Value = 0.017
$ws.Range("W7").Value = 0.284
$ws.Range("AI7").Value = 0.309
$ws.Range("AU7").Value = 0.317
$ws.Range("BA7").Value = 1.423
$ws.Range("BG7").Value = 0.432
$ws.Range("BM7").Value = 0.578
$ws.Range("BP7").Value = 0.474
$ws.Range("BQ7").Value = 0.488
# Row 8
$ws.Range("E8").Value = 0.607
$ws.Range("F8").Value = 0.11
$ws.Range("G8").Value = 0.332
$ws.Range("N8").Value = 0.779
$ws.Range("O8").Value = 0.061
$ws.Range("P8").Value = 0.246
$ws.Range("Q8").Value = 0.018
$ws.Range("W8").Value = 0.311
$ws.Range("AI8").Value = 0.332
$ws.Range("AJ8").Value = 0.128
$ws.Range("AK8").Value = 0.358
$ws.Range("AU8").Value = 0.318
$ws.Range("AW8").Value = 0.29
$ws.Range("BA8").Value = 1.732
$ws.Range("BB8").Value = 0.125
$ws.Range("BC8").Value = 0.353
$ws.Range("BG8").Value = 0.5590000000000001
$ws.Range("BH8").Value = 0.107
$ws.Range("BI8").Value = 0.328
$ws.Range("BM8").Value = 0.6919999999999999
$ws.Range("BN8").Value = 0.067
$ws.Range("BO8").Value = 0.258
$ws.Range("BP8").Value = 0.577
$ws.Range("BQ8").Value = 0.601
# Row 9
$ws.Range("E9").Value = 0.549
$ws.Range("F9").Value = 0.248
$ws.Range("G9").Value = 0.498
$ws.Range("N9").Value = 0.681
$ws.Range("O9").Value = 0.217
$ws.Range("P9").Value = 0.466
$ws.Range("W9").Value = 0.209
$ws.Range("X9").Value = 0.165
$ws.Range("Y9").Value = 0.406
$ws.Range("AI9").Value = 0.253
$ws.Range("AJ9").Value = 0.189
$ws.Range("AK9").Value = 0.435
$ws.Range("BA9").Value = 1.681
$ws.Range("BB9").Value = 0.246
$ws.Range("BC9").Value = 0.496
$ws.Range("BG9").Value = 0.593
$ws.Range("BH9").Value = 0.241
$ws.Range("BI9").Value = 0.491
$ws.Range("BM9").Value = 0.648
$ws.Range("BN9").Value = 0.228
$ws.Range("BO9").Value = 0.477
$ws.Range("BP9").Value = 0.5600000000000001
$ws.Range("BQ9").Value = 0.58
# Row 10
$ws.Range("E10").Value = 0.681
$ws.Range("F10").Value = 0.217
$ws.Range("G10").Value = 0.466
$ws.Range("N10").Value = 0.879
$ws.Range("O10").Value = 0.106
$ws.Range("P10").Value = 0.326
$ws.Range("W10").Value = 0.385
$ws.Range("X10").Value = 0.237
$ws.Range("Y10").Value = 0.487
$ws.Range("AI10").Value = 0.363
$ws.Range("AJ10").Value = 0.231
$ws.Range("AK10").Value = 0.481
$ws.Range("AU10").Value = 0.308
$ws.Range("AV10").Value = 0.213
$ws.Range("AW10").Value = 0.462
$ws.Range("BA10").Value = 2.065
$ws.Range("BB10").Value = 0.245
$ws.Range("BC10").Value = 0.495
$ws.Range("BG10").Value = 0.648
$ws.Range("BH10").Value = 0.228
$ws.Range("BI10").Value = 0.477
$ws.Range("BM10").Value = 0.846
$ws.Range("BN10").Value = 0.13
$ws.Range("BO10").Value = 0.361
$ws.Range("BP10").Value = 0.6879999999999999
$ws.Range("BQ10").Value = 0.721
# Row 11
$ws.Range("E11").Value = 0.714
$ws.Range("F11").Value = 0.204
$ws.Range("G11").Value = 0.452
$ws.Range("N11").Value = 0.901
$ws.Range("O11").Value = 0.089
$ws.Range("P11").Value = 0.299
$ws.Range("W11").Value = 0.385
$ws.Range("X11").Value = 0.237
$ws.Range("Y11").Value = 0.487
$ws.Range("AI11").Value = 0.396
$ws.Range("AJ11").Value = 0.239
$ws.Range("AK11").Value = 0.489
$ws.Range("AU11").Value = 0.451
$ws.Range("AV11").Value = 0.248
$ws.Range("AW11").Value = 0.498
$ws.Range("BA11").Value = 2.065
$ws.Range("BB11").Value = 0.245
$ws.Range("BC11").Value = 0.495
$ws.Range("BG11").Value = 0.648
$ws.Range("BH11").Value = 0.228
$ws.Range("BI11").Value = 0.477
$ws.Range("BM11").Value = 0.846
$ws.Range("BN11").Value = 0.13
$ws.Range("BO11").Value = 0.361
$ws.Range("BP11").Value = 0.6879999999999999
$ws.Range("BQ11").Value = 0.724
# Row 12
$ws.Range("E12").Value = 1.415
$ws.Range("F12").Value = 0.766
$ws.Range("G12").Value = 0.875
$ws.Range("N12").Value = 1.476
$ws.Range("O12").Value = 1.059
$ws.Range("P12").Value = 1.029
$ws.Range("W12").Value = 1.629
$ws.Range("X12").Value = 0.576
$ws.Range("Y12").Value = 0.759
$ws.Range("AI12").Value = 1.722
$ws.Range("AJ12").Value = 1.312
$ws.Range("AK12").Value = 1.145
$ws.Range("AU12").Value = 2.767
$ws.Range("AV12").Value = 2.737
$ws.Range("AW12").Value = 1.654
$ws.Range("BA12").Value = 3.728
$ws.Range("BB12").Value = 0.412
$ws.Range("BC12").Value = 0.642
$ws.Range("BG12").Value = 1.102
$ws.Range("BH12").Value = 0.125
$ws.Range("BI12").Value = 0.354
$ws.Range("BM12").Value = 1.299
$ws.Range("BN12").Value = 0.339
$ws.Range("BO12").Value = 0.583
$ws.Range("BP12").Value = 1.243
$ws.Range("BQ12").Value = 1.266
# Row 13
$ws.Range("E13").Value = 1.58
$ws.Range("F13").Value = 0.664
$ws.Range("G13").Value = 0.8149999999999999
$ws.Range("N13").Value = 2.065
$ws.Range("O13").Value = 0.957
$ws.Range("P13").Value = 0.978
$ws.Range("W13").Value = 1.031
$ws.Range("X13").Value = 0.193
$ws.Range("Y13").Value = 0.439
$ws.Range("AI13").Value = 1.284
$ws.Range("AJ13").Value = 0.374
$ws.Range("AK13").Value = 0.611
$ws.Range("AU13").Value = 2.285
$ws.Range("AV13").Value = 0.925
$ws.Range("AW13").Value = 0.962
$ws.Range("BA13").Value = 2.37
$ws.Range("BB13").Value = 0.302
$ws.Range("BC13").Value = 0.549
$ws.Range("BG13").Value = 0.585
$ws.Range("BH13").Value = 0.073
$ws.Range("BI13").Value = 0.27
$ws.Range("BM13").Value = 0.906
$ws.Range("BN13").Value = 0.284
$ws.Range("BO13").Value = 0.533
$ws.Range("BP13").Value = 0.79
$ws.Range("BQ13").Value = 0.73
